$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the style from the last existing date cell (A19) down to the new
# row's date cell (A20) so it keeps the same date number format/border/bold.
$ws.Range("A19").Copy($ws.Range("A20"))

# Populate the new row of data (row 20)
$ws.Range("A20").Value = 45986
$ws.Range("B20").Value = 2025
$ws.Range("C20").Value = 2.560577522109297
$ws.Range("D20").Value = 2026
$ws.Range("E20").Value = 1.676143333484292
